$d = $word.ActiveDocument
$d.TrackRevisions = $false

function Get-ParagraphRangeContaining([string]$needle) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p.Range
        }
    }
    throw ("Paragraph containing '" + $needle + "' not found")
}

# --- Part 1: "Palabras clave del recurso" keyword paragraph -----------------
# Split "Moléculas, estados de la materia, sólido, líquido, gaseoso." into
# several runs wrapped with proofErr spellStart/spellEnd markers, drop the
# spaces between the comma separated keywords, and relocate the stray
# "_GoBack" bookmark from the other paragraph in here (right after "eoso").

$para1 = Get-ParagraphRangeContaining("Palabras clave del recurso")

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="36FAF20C" w14:textId="48282049" w:rsidR="00CD652E" w:rsidRPr="006D02A8" w:rsidRDefault="00CD652E" w:rsidP="00CD652E"><w:pPr><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r w:rsidRPr="006D02A8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:b/><w:color w:val="FF0000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>*</w:t></w:r><w:r w:rsidRPr="006D02A8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidRPr="006D02A8"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:highlight w:val="green"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>Palabras clave del recurso (separadas por comas ",")</w:t></w:r><w:r w:rsidR="00B87883"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>m</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>oléculas,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>estados</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve"> de la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>materia,sólido,líquido,</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>gas</w:t></w:r><w:r w:rsidR="00EA3242"><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>eoso</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para1.InsertXML($xml1)

# --- Part 2: "ofrecen muy poca oposición" paragraph -------------------------
# The old stray "_GoBack" bookmark used to sit between the single leading
# space and "muy poca"; it moved to the paragraph above, so here the space
# run and the "muy poca" run simply collapse into one run.

$para2 = Get-ParagraphRangeContaining("oposición a cambiar su forma y volumen")

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="0723E241" w14:textId="3AA0B061" w:rsidR="009C5E8E" w:rsidRPr="009A13A2" w:rsidRDefault="009C5E8E" w:rsidP="009C5E8E"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl"/></w:rPr></w:pPr><w:r w:rsidRPr="009A13A2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t xml:space="preserve">En </w:t></w:r><w:r w:rsidRPr="009A13A2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>estado gaseoso</w:t></w:r><w:r w:rsidRPr="009A13A2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl"/></w:rPr><w:t>, las moléculas no tienen unión y sí mucha movilidad</w:t></w:r><w:r w:rsidR="0026726E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl" w:eastAsia="es-ES"/></w:rPr><w:t>;</w:t></w:r><w:r w:rsidR="00DD1B26"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl" w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="008456FA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl" w:eastAsia="es-ES"/></w:rPr><w:t>ofrecen</w:t></w:r><w:r w:rsidR="00DD1B26"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl" w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> muy poca</w:t></w:r><w:r w:rsidRPr="009A13A2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:lang w:val="es-ES_tradnl" w:eastAsia="es-ES"/></w:rPr><w:t xml:space="preserve"> oposición a cambiar su forma y volumen.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$para2.InsertXML($xml2)
